# Update the "word -> code" token table on Sheet1.
#  - Row 19 (B/C) changes from "2"/80 to "a"/81.
#  - Rows 33-44 (B/C) are rewritten with a shifted token stream.
#  - Rows 45-59 are brand-new rows continuing that token stream
#    (column A keeps incrementing, using the same bold/centered/
#    bordered style already used by the rest of column A).
#
# Two quirks of this engine's Range.Value setter need workarounds so the
# "word" column keeps storing plain text (matching every existing row):
#   1. A leading "=" is parsed as the start of a formula (like typing it
#      into a live grid), so it can't be poked in with a plain string
#      assignment.
#   2. A value that looks like a bare integer ("1", "2", "4", ...) is
#      auto-coerced to a Number cell instead of staying Text.
# The fix for both is the same: copy an existing cell that already holds
# the literal text we want and PasteSpecial just the value into the
# destination, instead of assigning through .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 0. Rescue the one token we're about to overwrite but still need -
# B37 currently holds the literal text "4", which is needed again at the
# new B39. Copy it across before row 37 itself gets overwritten below.
$ws.Range("B37").Copy()
$ws.Range("B39").PasteSpecial(-4163)

# ---- 1. Extend column A's styled range down to row 59 --------------
# Copy the formatting (bold font, thin border, centered/top alignment)
# already used by A2:A44 onto the new A45:A59 cells.
$ws.Range("A44").Copy()
for ($r = 45; $r -le 59; $r++) {
    $ws.Range("A$r").PasteSpecial(-4122)
}

# ---- 2. Column A index values for the new rows ----------------------
$aValues = @{
    45 = 43
    46 = 44
    47 = 45
    48 = 46
    49 = 47
    50 = 48
    51 = 49
    52 = 50
    53 = 51
    54 = 52
    55 = 53
    56 = 54
    57 = 55
    58 = 56
    59 = 57
}
foreach ($r in $aValues.Keys) {
    $ws.Range("A$r").Value = $aValues[$r]
}

# ---- 3. Target token ("word") / code pairs for columns B and C ------
$tokens = @{
    19 = @("a", 81)
    33 = @("+", 35)
    34 = @("1", 80)
    35 = @("*", 37)
    36 = @("(", 61)
    37 = @("2", 80)
    38 = @("+", 35)
    39 = @("4", 80)
    40 = @(")", 62)
    41 = @("e", 81)
    42 = @("=", 46)
    43 = @("1", 80)
    44 = @("+", 35)
    45 = @("2", 80)
    46 = @(")", 62)
    47 = @("f", 81)
    48 = @("=", 46)
    49 = @("(", 61)
    50 = @("1", 80)
    51 = @("*", 37)
    52 = @("1", 80)
    53 = @("*", 37)
    54 = @("(", 61)
    55 = @("1", 80)
    56 = @("*", 37)
    57 = @("2", 80)
    58 = @(")", 62)
    59 = @("#", 43)
}

# Stable cells (never touched by this edit) that already hold the exact
# literal text of the "numeric-looking" tokens / the "=" token, used as
# PasteSpecial(values) sources so the destination stays a Text cell.
$textSource = @{
    "1" = "B5"
    "2" = "B12"
    "=" = "B3"
}

foreach ($r in ($tokens.Keys | Sort-Object)) {
    $word = $tokens[$r][0]
    $code = $tokens[$r][1]

    if ($r -eq 39) {
        # already populated in step 0 above (rescued from the old B37)
    } elseif ($textSource.ContainsKey($word)) {
        $ws.Range($textSource[$word]).Copy()
        $ws.Range("B$r").PasteSpecial(-4163)
    } else {
        $ws.Range("B$r").Value = $word
    }

    $ws.Range("C$r").Value = $code
}
